$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet 1"

# Clear old content
$ws.Range("C3").Value = $null

# Populate A2:A7 with JSON strings
$ws.Range("A2").Value = '{"name":"Dreams & Nightmares Case","price":"$1.33 USD"}'
$ws.Range("A3").Value = '{"name":"Dreams & Nightmares Case","price":"$1.36 USD"}'
$ws.Range("A4").Value = '{"name":"Dreams & Nightmares Case","price":"$1.37 USD"}'
$ws.Range("A5").Value = '{"name":"Dreams & Nightmares Case","price":"$1.37 USD"}'
$ws.Range("A6").Value = '{"name":"Dreams & Nightmares Case","price":"$1.37 USD"}'
$ws.Range("A7").Value = '{"name":"Dreams & Nightmares Case","price":"$1.38 USD"}'
